# Rerun models and create results figures and tables
# Rename all worksheets (sheets keep their position/sheetId/order,
# only the sheet "name" attribute changes to new summ<number> identifiers).
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ51676947"
$wb.Worksheets.Item(2).Name = "summ51946245"
$wb.Worksheets.Item(3).Name = "summ52228963"
$wb.Worksheets.Item(4).Name = "summ52496427"
$wb.Worksheets.Item(5).Name = "summ52780550"
$wb.Worksheets.Item(6).Name = "summ53122407"
$wb.Worksheets.Item(7).Name = "summ53457002"
$wb.Worksheets.Item(8).Name = "summ53750928"
$wb.Worksheets.Item(9).Name = "summ54094219"
$wb.Worksheets.Item(10).Name = "summ54491583"
$wb.Worksheets.Item(11).Name = "summ54835153"
$wb.Worksheets.Item(12).Name = "summ55178977"
$wb.Worksheets.Item(13).Name = "summ55534395"
$wb.Worksheets.Item(14).Name = "summ55906614"
$wb.Worksheets.Item(15).Name = "summ56249239"
$wb.Worksheets.Item(16).Name = "summ56654351"
$wb.Worksheets.Item(17).Name = "summ57005969"
$wb.Worksheets.Item(18).Name = "summ57339664"
$wb.Worksheets.Item(19).Name = "summ57680018"
$wb.Worksheets.Item(20).Name = "summ58040755"
$wb.Worksheets.Item(21).Name = "summ58430328"
$wb.Worksheets.Item(22).Name = "summ58768172"
$wb.Worksheets.Item(23).Name = "summ59118332"
$wb.Worksheets.Item(24).Name = "summ59474780"
$wb.Worksheets.Item(25).Name = "summ00090195"
$wb.Worksheets.Item(26).Name = "summ00417975"
$wb.Worksheets.Item(27).Name = "summ00794533"
$wb.Worksheets.Item(28).Name = "summ01161972"
$wb.Worksheets.Item(29).Name = "summ01549354"
$wb.Worksheets.Item(30).Name = "summ01887530"
$wb.Worksheets.Item(31).Name = "summ02222027"
$wb.Worksheets.Item(32).Name = "summ02565763"
$wb.Worksheets.Item(33).Name = "summ02944045"
$wb.Worksheets.Item(34).Name = "summ03259793"
$wb.Worksheets.Item(35).Name = "summ03599206"
$wb.Worksheets.Item(36).Name = "summ03961701"
$wb.Worksheets.Item(37).Name = "summ04314453"
$wb.Worksheets.Item(38).Name = "summ04692446"
$wb.Worksheets.Item(39).Name = "summ05028913"
$wb.Worksheets.Item(40).Name = "summ05383298"
$wb.Worksheets.Item(41).Name = "summ05784860"
$wb.Worksheets.Item(42).Name = "summ06148439"
$wb.Worksheets.Item(43).Name = "summ06542339"
$wb.Worksheets.Item(44).Name = "summ06865440"
$wb.Worksheets.Item(45).Name = "summ07247659"
$wb.Worksheets.Item(46).Name = "summ07597408"
$wb.Worksheets.Item(47).Name = "summ07930766"
$wb.Worksheets.Item(48).Name = "summ08268510"
$wb.Worksheets.Item(49).Name = "summ08610598"
$wb.Worksheets.Item(50).Name = "summ08957931"
